# Update crypto price/volume data per the Sep 11 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.274.16"
Set-TextValue "E2" "  -2.94%  "
Set-TextValue "D3" "1.552.62"
Set-TextValue "E3" "  -4.78%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "207.11"
Set-TextValue "E5" "  -3.38%  "
Set-TextValue "E6" "  +0.00%  "
Set-TextValue "E7" "  -5.07%  "
Set-TextValue "D8" "0.0610"
Set-TextValue "E8" "  -1.56%  "
Set-TextValue "D9" "0.243"
Set-TextValue "E9" "  -3.24%  "
Set-TextValue "D10" "17.72"
Set-TextValue "E10" "  -4.81%  "
Set-TextValue "D11" "0.0781"
Set-TextValue "E11" "  -1.04%  "
Set-TextValue "D12" "1.767.04"
Set-TextValue "E12" "  -4.84%  "
Set-TextValue "D13" "1.549.69"
Set-TextValue "E13" "  -5.05%  "
Set-TextValue "D14" "3.99"
Set-TextValue "E14" "  -4.65%  "
Set-TextValue "D15" "0.505"
Set-TextValue "E15" "  -4.57%  "
Set-TextValue "D16" "25.257.21"
Set-TextValue "D17" "0.0₃0708"
Set-TextValue "E17" "  -4.61%  "
Set-TextValue "D18" "58.69"
Set-TextValue "E18" "  -4.76%  "
Set-TextValue "D19" "1.01"
Set-TextValue "D20" "185.60"
Set-TextValue "E20" "  -3.90%  "
Set-TextValue "D21" "4.11"
Set-TextValue "E21" "  -3.72%  "
Set-TextValue "D22" "9.27"
Set-TextValue "E22" "  -3.06%  "
Set-TextValue "D23" "5.85"
Set-TextValue "E23" "  -3.85%  "
Set-TextValue "E24" "  -4.36%  "
Set-TextValue "E25" "  -0.04%  "
Set-TextValue "D26" "139.69"
Set-TextValue "E26" "  -3.17%  "
Set-TextValue "E27" "  -5.19%  "
Set-TextValue "D28" "14.85"
Set-TextValue "D29" "6.39"
Set-TextValue "E29" "  -5.18%  "
Set-TextValue "E30" "  -6.87%  "
Set-TextValue "D31" "0.0466"
Set-TextValue "E31" "  -3.50%  "
Set-TextValue "D32" "3.03"
Set-TextValue "E32" "  -3.47%  "
Set-TextValue "E33" "  -5.23%  "
Set-TextValue "E34" "  -2.97%  "
Set-TextValue "D35" "2.33"
Set-TextValue "E35" "  -3.41%  "
Set-TextValue "D36" "1.084.70"
Set-TextValue "E36" "  -3.55%  "
Set-TextValue "E37" "  -0.09%  "
Set-TextValue "E38" "  -3.07%  "
Set-TextValue "D39" "0.495"
Set-TextValue "E39" "  -5.03%  "
Set-TextValue "B40" "ARBITRUM"
Set-TextValue "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.765"
Set-TextValue "E40" "  -10.45%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.24"
Set-TextValue "E41" "  -7.78%  "
Set-TextValue "D42" "0.801"
Set-TextValue "E42" "  +5.91%  "
Set-TextValue "D43" "92.54"
Set-TextValue "E43" "  -5.84%  "
Set-TextValue "D44" "5.04"
Set-TextValue "E44" "  -1.89%  "
Set-TextValue "D45" "1.681.91"
Set-TextValue "E45" "  -4.80%  "
Set-TextValue "E46" "  -2.44%  "
Set-TextValue "E47" "  -1.93%  "
Set-TextValue "D48" "52.33"
Set-TextValue "E48" "  -4.19%  "
Set-TextValue "D49" "0.0503"
Set-TextValue "E49" "  -3.99%  "
Set-TextValue "E50" "  -0.18%  "
Set-TextValue "E51" "  -2.22%  "
